$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.049.90'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '3.419.66'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.08'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.58%  '
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '4.001.57'
$ws.Range("E12").Value = '  -1.88%  '
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.57%  '
$ws.Range("D15").Value = '3.423.99'
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").Value = '62.006.28'
$ws.Range("E17").Value = '  -2.35%  '
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.565'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '3.554.62'
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000112'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.180'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.90'
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = '  -1.27%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.68%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.22%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '168.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").Value = '3.450.42'
$ws.Range("E40").Value = '  -1.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0777'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("E45").Value = '  -3.81%  '
$ws.Range("E46").Value = '  -5.48%  '
$ws.Range("D47").Value = '2.545.25'
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.18'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
